$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 265072.34
$ws.Range("I132").Value = 380745.4
$ws.Range("J132").Value = 18303.133
$ws.Range("K132").Value = 1142236.2
$ws.Range("L132").Value = 54909.399
$ws.Range("M132").Value = -1139706.2
$ws.Range("N132").Value = -59969.399
$ws.Range("H137").Value = 41668324
$ws.Range("I137").Value = 58824584
$ws.Range("K137").Value = 176473752
$ws.Range("M137").Value = -176471202
$ws.Range("H138").Value = 7409496
$ws.Range("J138").Value = 10641104
$ws.Range("L138").Value = 31923312
$ws.Range("N138").Value = -31933592

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3784.5715
$ws.Range("I32").Value = 1850.449
$ws.Range("J32").Value = 10554
$ws.Range("K32").Value = 1850.449
$ws.Range("L32").Value = 10554
$ws.Range("M32").Value = -1563.449
$ws.Range("N32").Value = -11128
$ws.Range("H61").Value = 4094.7273
$ws.Range("I61").Value = 2939
$ws.Range("K61").Value = 2939
$ws.Range("M61").Value = -2727
$ws.Range("H74").Value = 8902.529
$ws.Range("I74").Value = 1983.4546
$ws.Range("J74").Value = 21587.5
$ws.Range("K74").Value = 1983.4546
$ws.Range("L74").Value = 21587.5
$ws.Range("M74").Value = -1109.4546
$ws.Range("N74").Value = -23335.5
$ws.Range("H77").Value = 8902.529
$ws.Range("I77").Value = 1983.4546
$ws.Range("J77").Value = 21587.5
$ws.Range("K77").Value = 9917.273000000001
$ws.Range("L77").Value = 107937.5
$ws.Range("M77").Value = -5549.273000000001
$ws.Range("N77").Value = -116673.5
$ws.Range("H88").Value = 4347.5
$ws.Range("I88").Value = 1726.6666
$ws.Range("J88").Value = 5920
$ws.Range("K88").Value = 1726.6666
$ws.Range("L88").Value = 5920
$ws.Range("M88").Value = -1320.6666
$ws.Range("N88").Value = -6732
$ws.Range("H91").Value = 4347.5
$ws.Range("I91").Value = 1726.6666
$ws.Range("J91").Value = 5920
$ws.Range("K91").Value = 1726.6666
$ws.Range("L91").Value = 5920
$ws.Range("M91").Value = -322.6666
$ws.Range("N91").Value = -8728
$ws.Range("H132").Value = 3432.1667
$ws.Range("I132").Value = 2935.6667
$ws.Range("J132").Value = 4176.9165
$ws.Range("K132").Value = 8807.000100000001
$ws.Range("L132").Value = 12530.7495
$ws.Range("M132").Value = -6277.000100000001
$ws.Range("N132").Value = -17590.7495
$ws.Range("H136").Value = 4094.7273
$ws.Range("I136").Value = 2939
$ws.Range("K136").Value = 8817
$ws.Range("M136").Value = -6267

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2125
$ws.Range("I99").Value = 1750
$ws.Range("J99").Value = 2500
$ws.Range("K99").Value = 1750
$ws.Range("L99").Value = 2500
$ws.Range("M99").Value = -252
$ws.Range("N99").Value = -5496
$ws.Range("H124").Value = 36995
$ws.Range("J124").Value = 36995
$ws.Range("L124").Value = 36995
$ws.Range("N124").Value = -46815
$ws.Range("H134").Value = 25644030
$ws.Range("I134").Value = 38463276
$ws.Range("K134").Value = 115389828
$ws.Range("M134").Value = -115387293

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1305.0968
$ws.Range("I31").Value = 1246.3448
$ws.Range("J31").Value = 2157
$ws.Range("K31").Value = 1246.3448
$ws.Range("L31").Value = 2157
$ws.Range("M31").Value = -951.3448000000001
$ws.Range("N31").Value = -2747
$ws.Range("H34").Value = 1305.0968
$ws.Range("I34").Value = 1246.3448
$ws.Range("J34").Value = 2157
$ws.Range("K34").Value = 1246.3448
$ws.Range("L34").Value = 2157
$ws.Range("M34").Value = -1044.3448
$ws.Range("N34").Value = -2561
$ws.Range("H58").Value = 2289.2593
$ws.Range("I58").Value = 1584.5
$ws.Range("J58").Value = 3698.7778
$ws.Range("K58").Value = 1584.5
$ws.Range("L58").Value = 3698.7778
$ws.Range("M58").Value = -1381.5
$ws.Range("N58").Value = -4104.7778
$ws.Range("H122").Value = 1784.4
$ws.Range("I122").Value = 938.75
$ws.Range("K122").Value = 2816.25
$ws.Range("M122").Value = -366.25
$ws.Range("H132").Value = 2337.8948
$ws.Range("I132").Value = 1617.6207
$ws.Range("J132").Value = 4658.778
$ws.Range("K132").Value = 4852.8621
$ws.Range("L132").Value = 13976.334
$ws.Range("M132").Value = -2322.8621
$ws.Range("N132").Value = -19036.334
$ws.Range("H134").Value = 2505.1428
$ws.Range("I134").Value = 1092.3334
$ws.Range("K134").Value = 3277.0002
$ws.Range("M134").Value = -742.0001999999999
$ws.Range("H136").Value = 2289.2593
$ws.Range("I136").Value = 1584.5
$ws.Range("J136").Value = 3698.7778
$ws.Range("K136").Value = 4753.5
$ws.Range("L136").Value = 11096.3334
$ws.Range("M136").Value = -2203.5
$ws.Range("N136").Value = -16196.3334

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H103").Value = 1801.08
$ws.Range("I103").Value = 540.63635
$ws.Range("J103").Value = 2791.4285
$ws.Range("K103").Value = 1621.90905
$ws.Range("L103").Value = 8374.2855
$ws.Range("M103").Value = -742.90905
$ws.Range("N103").Value = -10132.2855
$ws.Range("H131").Value = 1528.7742
$ws.Range("I131").Value = 691.8
$ws.Range("J131").Value = 1602.193
$ws.Range("K131").Value = 2075.4
$ws.Range("L131").Value = 4806.579
$ws.Range("M131").Value = 2964.6
$ws.Range("N131").Value = -14886.579
$ws.Range("H132").Value = 1237.9231
$ws.Range("I132").Value = 304
$ws.Range("J132").Value = 1518.1
$ws.Range("K132").Value = 2736
$ws.Range("L132").Value = 13662.9
$ws.Range("M132").Value = -206
$ws.Range("N132").Value = -18722.9

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1508.5
$ws.Range("I102").Value = 1281.6666
$ws.Range("J102").Value = 1644.6
$ws.Range("K102").Value = 1281.6666
$ws.Range("L102").Value = 1644.6
$ws.Range("M102").Value = 340.3334
$ws.Range("N102").Value = -4888.6
$ws.Range("H132").Value = 3065.8857
$ws.Range("I132").Value = 2855.1428
$ws.Range("J132").Value = 3382
$ws.Range("K132").Value = 8565.428400000001
$ws.Range("L132").Value = 10146
$ws.Range("M132").Value = -6035.428400000001
$ws.Range("N132").Value = -15206
$ws.Range("H134").Value = 28406.5
$ws.Range("J134").Value = 28406.5
$ws.Range("L134").Value = 85219.5
$ws.Range("N134").Value = -90289.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 11269
$ws.Range("I22").Value = 1278
$ws.Range("J22").Value = 21260
$ws.Range("K22").Value = 1278
$ws.Range("L22").Value = 21260
$ws.Range("M22").Value = -983
$ws.Range("N22").Value = -21850
$ws.Range("H27").Value = 11269
$ws.Range("I27").Value = 1278
$ws.Range("J27").Value = 21260
$ws.Range("K27").Value = 1278
$ws.Range("L27").Value = 21260
$ws.Range("M27").Value = -1171
$ws.Range("N27").Value = -21474
$ws.Range("H46").Value = 905.38464
$ws.Range("I46").Value = 825.125
$ws.Range("J46").Value = 1033.8
$ws.Range("K46").Value = 825.125
$ws.Range("L46").Value = 1033.8
$ws.Range("M46").Value = -637.125
$ws.Range("N46").Value = -1409.8
$ws.Range("H132").Value = 5792.8335
$ws.Range("I132").Value = 4821.6
$ws.Range("K132").Value = 14464.8
$ws.Range("M132").Value = -11934.8
$ws.Range("H136").Value = 4640.2905
$ws.Range("I136").Value = 2608.125
$ws.Range("J136").Value = 6807.933
$ws.Range("K136").Value = 7824.375
$ws.Range("L136").Value = 20423.799
$ws.Range("M136").Value = -5274.375
$ws.Range("N136").Value = -25523.799

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1642.65
$ws.Range("I122").Value = 1439.0769
$ws.Range("K122").Value = 4317.2307
$ws.Range("M122").Value = -1867.2307
$ws.Range("H132").Value = 10640853
$ws.Range("I132").Value = 13890974
$ws.Range("J132").Value = 4094.0908
$ws.Range("K132").Value = 41672922
$ws.Range("L132").Value = 12282.2724
$ws.Range("M132").Value = -41670392
$ws.Range("N132").Value = -17342.2724
$ws.Range("H136").Value = 13375606
$ws.Range("I136").Value = 33434684
$ws.Range("J136").Value = 2887.0667
$ws.Range("K136").Value = 100304052
$ws.Range("L136").Value = 8661.2001
$ws.Range("M136").Value = -100301502
$ws.Range("N136").Value = -13761.2001
